# "update on categories and faculty"
#
# 1. Fix a typo in the Faculty/category shared string used throughout the
#    "Faculties & Programs" sheet: "Varcity Clubs" -> "Varsity Clubs".
# 2. Switch the active/selected sheet from "Faculty & Academic" to
#    "Faculties & Programs", updating each sheet's remembered selection.

$wb = $excel.ActiveWorkbook

$wsFaculties = $wb.Worksheets.Item("Faculties & Programs")
$wsFaculty   = $wb.Worksheets.Item("Faculty & Academic")

# --- 1. Correct the category typo everywhere it appears ---
$wsFaculties.Cells.Replace("Varcity Clubs", "Varsity Clubs")

# --- 2. Make "Faculties & Programs" the active tab/selection ---
$wsFaculty.Range("C241").Select()

$wsFaculties.Activate()
$wsFaculties.Range("E240").Select()
